# Updated distance sorting and flow on
# The team.name column (A) is re-ordered within each 18-row state block.
# Distances / states per row-position stay exactly where they are; only the
# team label shown at each row-position changes, matching the new sort order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newOrder = @(
    "Carlton",
    "Adelaide Crows",
    "GWS Giants",
    "Fremantle",
    "Gold Coast Suns",
    "Essendon",
    "Collingwood",
    "Brisbane Lions",
    "Melbourne",
    "Geelong Cats",
    "Hawthorn",
    "Port Adelaide",
    "North Melbourne",
    "St Kilda",
    "Richmond",
    "West Coast Eagles",
    "Sydney Swans",
    "Western Bulldogs"
)

$blockStarts = @(2, 20, 38, 56, 74, 92, 110, 128, 146)

foreach ($start in $blockStarts) {
    for ($i = 0; $i -lt $newOrder.Length; $i++) {
        $row = $start + $i
        $ws.Cells.Item($row, 1).Value = $newOrder[$i]
    }
}
